$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.985.43"
$ws.Range("E2").Value = "  +2.80%  "

$ws.Range("D3").Value = "1.599.46"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.21"
$ws.Range("E5").Value = "  +2.35%  "

$ws.Range("E7").Value = "  +1.43%  "

$ws.Range("E8").Value = "  +1.51%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("E10").Value = "  +1.79%  "

$ws.Range("E11").Value = "  +3.59%  "

$ws.Range("D12").Value = "1.823.23"
$ws.Range("E12").Value = "  +2.56%  "

$ws.Range("D13").Value = "1.600.04"
$ws.Range("E13").Value = "  +2.49%  "

$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").Value = "25.991.07"
$ws.Range("E16").Value = "  +2.76%  "

$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.52"
$ws.Range("E20").Value = "  +7.87%  "

$ws.Range("E21").Value = "  +2.70%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("E23").Value = "  +2.54%  "

$ws.Range("E24").Value = "  +7.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.55"
$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -7.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.13"
$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("E29").Value = "  +1.06%  "

$ws.Range("E30").Value = "  +1.84%  "

$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("E32").Value = "  +1.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("E35").Value = "  +2.39%  "

$ws.Range("D36").Value = "1.123.95"

$ws.Range("E37").Value = "  +10.42%  "

$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("E39").Value = "  -1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.788"
$ws.Range("E40").Value = "  +2.04%  "

$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "1.734.75"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.83"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.50"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.47"

$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.21"
$ws.Range("E51").Value = "  +0.17%  "
